$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 50, shifting existing rows 50-60 down to 51-61.
$ws.Rows.Item(50).Insert()

# Populate the newly inserted row 50 with the data synced from Markdown.
# Prefix the "year" column with an apostrophe so the numeric-looking text
# "2006" is stored as text (matches the other year cells in this column,
# which are all text such as "2006", "2004-2007", "nan", ...) instead of
# being auto-converted to a number.
$ws.Cells.Item(50, 1).Value = "'2006"
$ws.Cells.Item(50, 2).Value = "**神戸大学** <br> [市販ミンチ肉 における黄色ブドウ球菌汚染調査と分離株の性状](https://www.jstage.jst.go.jp/article/jsfm1994/23/4/23_4_217/_pdf/-char/ja)"
$ws.Cells.Item(50, 3).Value = "未登録"
